$wb = $excel.ActiveWorkbook

# Sheet ALC (index 1) Row 9
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 86  # H9: was 38.5
$ws.Cells.Item(9, 9).Value = 115  # I9: was 30
$ws.Cells.Item(9, 10).Value = 66.666664  # J9: was 55.5
$ws.Cells.Item(9, 11).Value = 115  # K9: was 30
$ws.Cells.Item(9, 12).Value = 66.666664  # L9: was 55.5
$ws.Cells.Item(9, 13).Value = 54  # M9: was 139
$ws.Cells.Item(9, 14).Value = -404.666664  # N9: was -393.5

# Sheet ALC (index 1) Row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 7936852.5  # H41: was 7407913.5
$ws.Cells.Item(41, 9).Value = 7936852.5  # I41: was 12346329
$ws.Cells.Item(41, 10).Value = 0  # J41: was 289.5
$ws.Cells.Item(41, 11).Value = 7936852.5  # K41: was 12346329
$ws.Cells.Item(41, 12).Value = 0  # L41: was 289.5
$ws.Cells.Item(41, 13).Value = -7936412.5  # M41: was -12345889
$ws.Cells.Item(41, 14).ClearContents()  # N41: was -1169.5

# Sheet ALC (index 1) Row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 258304.94  # H132: was 258318.38
$ws.Cells.Item(132, 9).Value = 265187.62  # I132: was 271072.8
$ws.Cells.Item(132, 10).Value = 100003  # J132: was 67002
$ws.Cells.Item(132, 11).Value = 795562.86  # K132: was 813218.3999999999
$ws.Cells.Item(132, 12).Value = 300009  # L132: was 201006
$ws.Cells.Item(132, 13).Value = -793032.86  # M132: was -810688.3999999999
$ws.Cells.Item(132, 14).Value = -305069  # N132: was -206066

# Sheet ALC (index 1) Row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 16667598  # H137: was 19608854
$ws.Cells.Item(137, 9).Value = 19608542  # I137: was 24390978
$ws.Cells.Item(137, 10).Value = 2253.111  # J137: was 2147.8
$ws.Cells.Item(137, 11).Value = 58825626  # K137: was 73172934
$ws.Cells.Item(137, 12).Value = 6759.333  # L137: was 6443.400000000001
$ws.Cells.Item(137, 13).Value = -58823076  # M137: was -73170384
$ws.Cells.Item(137, 14).Value = -11859.333  # N137: was -11543.4

# Sheet ARM (index 2) Row 46
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(46, 8).Value = 6727.75  # H46: was 7610.4
$ws.Cells.Item(46, 10).Value = 6727.75  # J46: was 7610.4
$ws.Cells.Item(46, 12).Value = 6727.75  # L46: was 7610.4
$ws.Cells.Item(46, 14).Value = -7365.75  # N46: was -8248.4

# Sheet ARM (index 2) Row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1689.7693  # H61: was 1719.8868
$ws.Cells.Item(61, 9).Value = 1368.8206  # I61: was 1392.2927
$ws.Cells.Item(61, 10).Value = 2652.6155  # J61: was 2839.1667
$ws.Cells.Item(61, 11).Value = 1368.8206  # K61: was 1392.2927
$ws.Cells.Item(61, 12).Value = 2652.6155  # L61: was 2839.1667
$ws.Cells.Item(61, 13).Value = -1156.8206  # M61: was -1180.2927
$ws.Cells.Item(61, 14).Value = -3076.6155  # N61: was -3263.1667

# Sheet ARM (index 2) Row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 2818.44  # H74: was 2766.4211
$ws.Cells.Item(74, 9).Value = 942.3090999999999  # I74: was 908.4828
$ws.Cells.Item(74, 10).Value = 7977.8  # J74: was 8753.111000000001
$ws.Cells.Item(74, 11).Value = 942.3090999999999  # K74: was 908.4828
$ws.Cells.Item(74, 12).Value = 7977.8  # L74: was 8753.111000000001
$ws.Cells.Item(74, 13).Value = -68.30909999999994  # M74: was -34.4828
$ws.Cells.Item(74, 14).Value = -9725.799999999999  # N74: was -10501.111

# Sheet ARM (index 2) Row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 2818.44  # H77: was 2766.4211
$ws.Cells.Item(77, 9).Value = 942.3090999999999  # I77: was 908.4828
$ws.Cells.Item(77, 10).Value = 7977.8  # J77: was 8753.111000000001
$ws.Cells.Item(77, 11).Value = 4711.5455  # K77: was 4542.414
$ws.Cells.Item(77, 12).Value = 39889  # L77: was 43765.55500000001
$ws.Cells.Item(77, 13).Value = -343.5455000000002  # M77: was -174.4139999999998
$ws.Cells.Item(77, 14).Value = -48625  # N77: was -52501.55500000001

# Sheet ARM (index 2) Row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 2932.625  # H122: was 2262.3333
$ws.Cells.Item(122, 9).Value = 2993.5  # I122: was 2482.625
$ws.Cells.Item(122, 10).Value = 2750  # J122: was 500
$ws.Cells.Item(122, 11).Value = 8980.5  # K122: was 7447.875
$ws.Cells.Item(122, 12).Value = 8250  # L122: was 1500
$ws.Cells.Item(122, 13).Value = -6530.5  # M122: was -4997.875
$ws.Cells.Item(122, 14).Value = -13150  # N122: was -6400

# Sheet ARM (index 2) Row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 2971.3333  # H132: was 2970
$ws.Cells.Item(132, 9).Value = 2174.4285  # I132: was 2172.7144
$ws.Cells.Item(132, 11).Value = 6523.2855  # K132: was 6518.1432
$ws.Cells.Item(132, 13).Value = -3993.2855  # M132: was -3988.1432

# Sheet ARM (index 2) Row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 1689.7693  # H136: was 1719.8868
$ws.Cells.Item(136, 9).Value = 1368.8206  # I136: was 1392.2927
$ws.Cells.Item(136, 10).Value = 2652.6155  # J136: was 2839.1667
$ws.Cells.Item(136, 11).Value = 4106.4618  # K136: was 4176.8781
$ws.Cells.Item(136, 12).Value = 7957.8465  # L136: was 8517.500100000001
$ws.Cells.Item(136, 13).Value = -1556.4618  # M136: was -1626.8781
$ws.Cells.Item(136, 14).Value = -13057.8465  # N136: was -13617.5001

# Sheet BSM (index 3) Row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 18520614  # H134: was 15626838
$ws.Cells.Item(134, 9).Value = 22728912  # I134: was 18869362
$ws.Cells.Item(134, 10).Value = 4102.6  # J134: was 3766
$ws.Cells.Item(134, 11).Value = 68186736  # K134: was 56608086
$ws.Cells.Item(134, 12).Value = 12307.8  # L134: was 11298
$ws.Cells.Item(134, 13).Value = -68184201  # M134: was -56605551
$ws.Cells.Item(134, 14).Value = -17377.8  # N134: was -16368

# Sheet CRP (index 4) Row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 50  # H7: was 49.333332
$ws.Cells.Item(7, 9).Value = 0  # I7: was 50
$ws.Cells.Item(7, 10).Value = 50  # J7: was 48.8
$ws.Cells.Item(7, 11).Value = 0  # K7: was 50
$ws.Cells.Item(7, 12).Value = 50  # L7: was 48.8
$ws.Cells.Item(7, 13).ClearContents()  # M7: was 63
$ws.Cells.Item(7, 14).Value = -276  # N7: was -274.8

# Sheet CRP (index 4) Row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1515.8096  # H31: was 1576.9
$ws.Cells.Item(31, 9).Value = 869.04  # I31: was 905.7447
$ws.Cells.Item(31, 11).Value = 869.04  # K31: was 905.7447
$ws.Cells.Item(31, 13).Value = -574.04  # M31: was -610.7447

# Sheet CRP (index 4) Row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 1515.8096  # H34: was 1576.9
$ws.Cells.Item(34, 9).Value = 869.04  # I34: was 905.7447
$ws.Cells.Item(34, 11).Value = 869.04  # K34: was 905.7447
$ws.Cells.Item(34, 13).Value = -667.04  # M34: was -703.7447

# Sheet CRP (index 4) Row 42
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(42, 8).Value = 5720  # H42: was 6050
$ws.Cells.Item(42, 9).Value = 4133.3335  # I42: was 5000
$ws.Cells.Item(42, 11).Value = 4133.3335  # K42: was 5000
$ws.Cells.Item(42, 13).Value = -3540.3335  # M42: was -4407

# Sheet CRP (index 4) Row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94, 8).Value = 567.85  # H94: was 682.5806
$ws.Cells.Item(94, 10).Value = 514.7692  # J94: was 695.8823
$ws.Cells.Item(94, 12).Value = 514.7692  # L94: was 695.8823
$ws.Cells.Item(94, 14).Value = -1416.7692  # N94: was -1597.8823

# Sheet CRP (index 4) Row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 2589.641  # H132: was 1934.3928
$ws.Cells.Item(132, 9).Value = 1967.3103  # I132: was 1470.186
$ws.Cells.Item(132, 10).Value = 4394.4  # J132: was 3469.8462
$ws.Cells.Item(132, 11).Value = 5901.9309  # K132: was 4410.558
$ws.Cells.Item(132, 12).Value = 13183.2  # L132: was 10409.5386
$ws.Cells.Item(132, 13).Value = -3371.9309  # M132: was -1880.558
$ws.Cells.Item(132, 14).Value = -18243.2  # N132: was -15469.5386

# Sheet CUL (index 5) Row 14
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 53.4  # H14: was 61.166668
$ws.Cells.Item(14, 9).Value = 53.4  # I14: was 61.166668
$ws.Cells.Item(14, 11).Value = 160.2  # K14: was 183.500004
$ws.Cells.Item(14, 13).Value = 12.80000000000001  # M14: was -10.50000399999999

# Sheet CUL (index 5) Row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 656.73334  # H92: was 789.44446
$ws.Cells.Item(92, 9).Value = 700.5  # I92: was 800.5
$ws.Cells.Item(92, 10).Value = 640.8182  # J92: was 780.6
$ws.Cells.Item(92, 11).Value = 2101.5  # K92: was 2401.5
$ws.Cells.Item(92, 12).Value = 1922.4546  # L92: was 2341.8
$ws.Cells.Item(92, 13).Value = -853.5  # M92: was -1153.5
$ws.Cells.Item(92, 14).Value = -4418.4546  # N92: was -4837.8

# Sheet CUL (index 5) Row 98
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(98, 8).Value = 388  # H98: was 385.85715
$ws.Cells.Item(98, 9).Value = 296.2857  # I98: was 175.25
$ws.Cells.Item(98, 10).Value = 548.5  # J98: was 666.6667
$ws.Cells.Item(98, 11).Value = 888.8571000000001  # K98: was 525.75
$ws.Cells.Item(98, 12).Value = 1645.5  # L98: was 2000.0001
$ws.Cells.Item(98, 13).Value = 609.1428999999999  # M98: was 972.25
$ws.Cells.Item(98, 14).Value = -4641.5  # N98: was -4996.0001

# Sheet CUL (index 5) Row 103
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(103, 8).Value = 431376.84  # H103: was 602127.6
$ws.Cells.Item(103, 9).Value = 2000  # I103: was 3000
$ws.Cells.Item(103, 10).Value = 603127.6  # J103: was 751909.5
$ws.Cells.Item(103, 11).Value = 6000  # K103: was 9000
$ws.Cells.Item(103, 12).Value = 1809382.8  # L103: was 2255728.5
$ws.Cells.Item(103, 13).Value = -5121  # M103: was -8121
$ws.Cells.Item(103, 14).Value = -1811140.8  # N103: was -2257486.5

# Sheet CUL (index 5) Row 110
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(110, 8).Value = 1133.3334  # H110: was 0
$ws.Cells.Item(110, 9).Value = 200  # I110: was 0
$ws.Cells.Item(110, 10).Value = 3000  # J110: was 0
$ws.Cells.Item(110, 11).Value = 600  # K110: was 0
$ws.Cells.Item(110, 12).Value = 9000  # L110: was 0
$ws.Cells.Item(110, 13).Value = 3490  # M110: was None
$ws.Cells.Item(110, 14).Value = -17180  # N110: was None

# Sheet CUL (index 5) Row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1589.375  # H131: was 1512.6571
$ws.Cells.Item(131, 10).Value = 2012.5  # J131: was 1866.037
$ws.Cells.Item(131, 12).Value = 6037.5  # L131: was 5598.111
$ws.Cells.Item(131, 14).Value = -16117.5  # N131: was -15678.111

# Sheet GSM (index 6) Row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2729.8103  # H132: was 2761.9824
$ws.Cells.Item(132, 9).Value = 2291.2888  # I132: was 2289.0222
$ws.Cells.Item(132, 10).Value = 4247.769  # J132: was 4535.5835
$ws.Cells.Item(132, 11).Value = 6873.866399999999  # K132: was 6867.0666
$ws.Cells.Item(132, 12).Value = 12743.307  # L132: was 13606.7505
$ws.Cells.Item(132, 13).Value = -4343.866399999999  # M132: was -4337.0666
$ws.Cells.Item(132, 14).Value = -17803.307  # N132: was -18666.7505

# Sheet LTW (index 7) Row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3553.1765  # H122: was 3544
$ws.Cells.Item(122, 9).Value = 2004  # I122: was 1902
$ws.Cells.Item(122, 10).Value = 3650  # J122: was 3778.5715
$ws.Cells.Item(122, 11).Value = 6012  # K122: was 5706
$ws.Cells.Item(122, 12).Value = 10950  # L122: was 11335.7145
$ws.Cells.Item(122, 13).Value = -3562  # M122: was -3256
$ws.Cells.Item(122, 14).Value = -15850  # N122: was -16235.7145

# Sheet LTW (index 7) Row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 4505.4814  # H132: was 4699.451
$ws.Cells.Item(132, 9).Value = 4592.3335  # I132: was 4852.6665
$ws.Cells.Item(132, 11).Value = 13777.0005  # K132: was 14557.9995
$ws.Cells.Item(132, 13).Value = -11247.0005  # M132: was -12027.9995

# Sheet LTW (index 7) Row 134
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(134, 8).Value = 0  # H134: was 37495
$ws.Cells.Item(134, 10).Value = 0  # J134: was 37495
$ws.Cells.Item(134, 12).Value = 0  # L134: was 37495
$ws.Cells.Item(134, 14).ClearContents()  # N134: was -47635

# Sheet LTW (index 7) Row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 3780.1177  # H136: was 3581.8704
$ws.Cells.Item(136, 9).Value = 2290.0732  # I136: was 2111.8445
$ws.Cells.Item(136, 10).Value = 9889.299999999999  # J136: was 10932
$ws.Cells.Item(136, 11).Value = 6870.219599999999  # K136: was 6335.5335
$ws.Cells.Item(136, 12).Value = 29667.9  # L136: was 32796
$ws.Cells.Item(136, 13).Value = -4320.219599999999  # M136: was -3785.5335
$ws.Cells.Item(136, 14).Value = -34767.89999999999  # N136: was -37896

# Sheet WVR (index 8) Row 108
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(108, 8).Value = 28800  # H108: was 23542
$ws.Cells.Item(108, 10).Value = 28800  # J108: was 23542
$ws.Cells.Item(108, 12).Value = 28800  # L108: was 23542
$ws.Cells.Item(108, 14).Value = -36480  # N108: was -31222

# Sheet WVR (index 8) Row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 6099001  # H132: was 5556893.5
$ws.Cells.Item(132, 9).Value = 8198047  # I132: was 7247583
$ws.Cells.Item(132, 11).Value = 24594141  # K132: was 21742749
$ws.Cells.Item(132, 13).Value = -24591611  # M132: was -21740219

# Sheet WVR (index 8) Row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 16381.738  # H136: was 14410.554
$ws.Cells.Item(136, 9).Value = 17839.797  # I136: was 15278.493
$ws.Cells.Item(136, 10).Value = 2044.1666  # J136: was 2433
$ws.Cells.Item(136, 11).Value = 53519.391  # K136: was 45835.479
$ws.Cells.Item(136, 12).Value = 6132.4998  # L136: was 7299
$ws.Cells.Item(136, 13).Value = -50969.391  # M136: was -43285.479
$ws.Cells.Item(136, 14).Value = -11232.4998  # N136: was -12399
